$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objetivos = 'Mostrar ao aluno o rio no seu funcionamento natural e a intervenção antrópica visando o uso da água. Mostrar ao aluno os usos dos recursos hídricos, suas demandas e suas limitações. Desenvolver junto ao aluno os procedimentos quantitativos relativos ás técnicas de formação de reservatórios e da implantação dos diversos usos da água, consumptivos ou não, apontando decorrentes impactos socioeconômicos e ambientais.'
$docentes = '7455355 - Robson da Silva Rocha'
$programaResumido = 'Barragens e Reservatórios. Usos da água demandados para o interesse humano e Panorama Geral da Engenharia dos Recursos Hídricos. Impactos Ambientais dos Usos da Água. Gestão dos Recursos Hídricos.'
$programa = '- Políticas Públicas, Balanço Hídrico,- Demanda de água e disponibilidade dos recursos hídricos: Abastecimento Humano, águas para Agropecuária e indústria. - Hidreletricidade. - Barragens e Reservatórios,- Navegação Interior.- Águas Subterrâneas.- Gerenciamento dos Recursos Hídricos.- Hidroeconomia- Relação entre saneamento e qualidade da água'
$metodo = 'Avaliação baseada em trabalhos com dados reais, exercícios, trabalhos práticos e relatórios.'
$criterio = 'Média ponderada das notas atribuídas aos exercícios e trabalhos práticos e relatórios.'
$norma = '1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0.'
$bibliografia = 'REBOUÇAS, A. C. et alli, Águas Doces do Brasil  Capital Ecológico, Usos e Conservação, 2a. ed. Escrituras Editora, São Paulo, 2002.GARCEZ, L. N. e ALVAREZ, G. A.  Hidrologia, Editora Edgard Blücher, ISBN 8521201699, 304 pgs, 2004.TUCCI, C.E.M. (organizador) - Hidrologia, Ciência e Aplicação, Coleção ABRH de Recursos Hídricos, vol. 4, EDUSP/ABRH, 1993: 943 p.VILLELA, S.M. & MATTOS, A. - Hidrologia Aplicada, McGraw-Hill do Brasil, 1975: 245 p.LINSLEY, R.K.; FRANZINI, J.B. - Engenharia de Recursos Hídricos, McGraw-Hill do Brasil, EDUSP, 1978: 793 p.DAKER, A. Hidráulica aplicada à agricultura. Livraria Freitas Bastos S.A., 1983, v.1.MOTA, S. (1995). Preservação e Conservação de Recursos Hídricos. ABES. sq. edição.FEITOSA, F. A. C. e MANOEL FILHO, J. - "Hidrogeologia - Conceitos e Aplicações, CPRM, Fortaleza, 391 p., 2000.SCHREIBER, G,.P. - Usinas Hidrelétricas - Editora Edgard Blucher, São Paulo.'

# Insert a new row at position 13 (shifts old rows 13-21 down to 14-22)
$ws.Rows.Item(13).Insert()

# Copy B/C formatting from the row below (the shifted former row 13, now row 14)
# onto the new row 13, so B13/C13 get the correct column styles (s=2 / s=3).
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# New row 13 only has B13/C13 populated with the teacher info; A13 stays empty.
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes

# Row 10 (Objetivos:) - replace the erroneous teacher text with the actual objectives text.
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Row 14 (Programa resumido:) - replace "Semestral" with the real summary text.
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 (Programa:) - replace the erroneous date with the real program text.
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 19 (Método:) - replace the erroneous teacher text with the real method text.
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 (Critério:) - replace with the real criteria text.
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 (Norma de recuperação:) - replace with the real recovery norm text.
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Row 22 (Bibliografia:) - replace with the real bibliography text.
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
